# Committing the changes for proper extent report
# The "Test Suite" sheet's Runmode column (C2:C7) should all read "Y"
# (previously a mix of "Y"/"N"), and the selection should reflect the
# full updated range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = "Y"

$ws.Range("C2:C7").Select()
